$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 135
$ws.Range("I9").Value = 36.5
$ws.Range("J9").Value = 266.33334
$ws.Range("K9").Value = 36.5
$ws.Range("L9").Value = 266.33334
$ws.Range("M9").Value = 132.5
$ws.Range("N9").Value = -604.33334
$ws.Range("H33").Value = 1493.6666
$ws.Range("I33").Value = 421.7143
$ws.Range("J33").Value = 2994.4
$ws.Range("K33").Value = 421.7143
$ws.Range("L33").Value = 2994.4
$ws.Range("M33").Value = -192.7143
$ws.Range("N33").Value = -3452.4
$ws.Range("H98").Value = 1643.5264
$ws.Range("I98").Value = 1398.1538
$ws.Range("K98").Value = 1398.1538
$ws.Range("M98").Value = 99.84619999999995
$ws.Range("H113").Value = 4211.3335
$ws.Range("I113").Value = 3974.9
$ws.Range("K113").Value = 3974.9
$ws.Range("M113").Value = -720.9000000000001
$ws.Range("H122").Value = 1643.5264
$ws.Range("I122").Value = 1398.1538
$ws.Range("K122").Value = 4194.4614
$ws.Range("M122").Value = -1744.4614
$ws.Range("H137").Value = 1876.4166
$ws.Range("I137").Value = 1592.4546
$ws.Range("K137").Value = 4777.3638
$ws.Range("M137").Value = -2227.3638

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1996.4445
$ws.Range("I45").Value = 1996.4445
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1996.4445
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1619.4445
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 10999.5
$ws.Range("I61").Value = 10999.5
$ws.Range("K61").Value = 10999.5
$ws.Range("M61").Value = -10787.5
$ws.Range("H88").Value = 3018.4167
$ws.Range("I88").Value = 604.4
$ws.Range("J88").Value = 4742.7144
$ws.Range("K88").Value = 604.4
$ws.Range("L88").Value = 4742.7144
$ws.Range("M88").Value = -198.4
$ws.Range("N88").Value = -5554.7144
$ws.Range("H91").Value = 3018.4167
$ws.Range("I91").Value = 604.4
$ws.Range("J91").Value = 4742.7144
$ws.Range("K91").Value = 604.4
$ws.Range("L91").Value = 4742.7144
$ws.Range("M91").Value = 799.6
$ws.Range("N91").Value = -7550.7144
$ws.Range("H136").Value = 10999.5
$ws.Range("I136").Value = 10999.5
$ws.Range("K136").Value = 32998.5
$ws.Range("M136").Value = -30448.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 30002
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 30002
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 30002
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -30348
$ws.Range("H105").Value = 1119.5
$ws.Range("J105").Value = 1449.5
$ws.Range("L105").Value = 1449.5
$ws.Range("N105").Value = -4943.5
$ws.Range("H134").Value = 10728.286
$ws.Range("J134").Value = 10000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -35070

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1029.75
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1029.75
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1029.75
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1729.75
$ws.Range("H31").Value = 3048.375
$ws.Range("I31").Value = 3048.375
$ws.Range("K31").Value = 3048.375
$ws.Range("M31").Value = -2753.375
$ws.Range("H34").Value = 3048.375
$ws.Range("I34").Value = 3048.375
$ws.Range("K34").Value = 3048.375
$ws.Range("M34").Value = -2846.375

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1000
$ws.Range("J34").Value = 1000
$ws.Range("L34").Value = 3000
$ws.Range("N34").Value = -3168
$ws.Range("H60").Value = 31.666666
$ws.Range("I60").Value = 31.666666
$ws.Range("K60").Value = 94.99999800000001
$ws.Range("M60").Value = 156.000002
$ws.Range("H114").Value = 279.5
$ws.Range("I114").Value = 228
$ws.Range("J114").Value = 331
$ws.Range("K114").Value = 684
$ws.Range("L114").Value = 993
$ws.Range("M114").Value = 2570
$ws.Range("N114").Value = -7501
$ws.Range("H121").Value = 810.8333
$ws.Range("J121").Value = 1087
$ws.Range("L121").Value = 3261
$ws.Range("N121").Value = -5881
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H139").Value = 994
$ws.Range("I139").Value = 994
$ws.Range("K139").Value = 2982
$ws.Range("M139").Value = 2158
$ws.Range("H140").Value = 3041.1333
$ws.Range("I140").Value = 2543.1667
$ws.Range("K140").Value = 7629.500100000001
$ws.Range("M140").Value = -2449.500100000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 425
$ws.Range("J13").Value = 600
$ws.Range("L13").Value = 600
$ws.Range("N13").Value = -878
$ws.Range("H70").Value = 7499.5
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730
$ws.Range("H73").Value = 7499.5
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064
$ws.Range("H80").Value = 2770
$ws.Range("I80").Value = 1486.2
$ws.Range("J80").Value = 3572.375
$ws.Range("K80").Value = 1486.2
$ws.Range("L80").Value = 3572.375
$ws.Range("M80").Value = -488.2
$ws.Range("N80").Value = -5568.375
$ws.Range("H83").Value = 2770
$ws.Range("I83").Value = 1486.2
$ws.Range("J83").Value = 3572.375
$ws.Range("K83").Value = 7431
$ws.Range("L83").Value = 17861.875
$ws.Range("M83").Value = -2439
$ws.Range("N83").Value = -27845.875
$ws.Range("H122").Value = 10420593
$ws.Range("I122").Value = 12503912
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 37511736
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -37509286
$ws.Range("N122").Value = -16900
$ws.Range("H126").Value = 3426.4285
$ws.Range("I126").Value = 2248.75
$ws.Range("J126").Value = 4996.6665
$ws.Range("K126").Value = 6746.25
$ws.Range("L126").Value = 14989.9995
$ws.Range("M126").Value = -4276.25
$ws.Range("N126").Value = -19929.9995
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3325
$ws.Range("I46").Value = 2414.4285
$ws.Range("K46").Value = 2414.4285
$ws.Range("M46").Value = -2226.4285
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 3247
$ws.Range("I122").Value = 2900.5
$ws.Range("K122").Value = 8701.5
$ws.Range("M122").Value = -6251.5
$ws.Range("H132").Value = 14371.25
$ws.Range("I132").Value = 17495
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 52485
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -49955
$ws.Range("N132").Value = -20060

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1989
$ws.Range("I132").Value = 1989
$ws.Range("K132").Value = 5967
$ws.Range("M132").Value = -3437
$ws.Range("H136").Value = 3355
$ws.Range("I136").Value = 3355
$ws.Range("K136").Value = 10065
$ws.Range("M136").Value = -7515
